# ADDED EMAIL ID IN POSTER
#
# Updates the "Acknowledgements" byline on the poster slide:
#   - "Priyanka Galla" -> "By: Priyanka Galla | Advisor: Dr. Nathan Eloe"
#   - "Northwest Missouri State University" -> "Email: s534884@nwmissouri.edu"
# and re-caches the master's datetimeFigureOut field text (7/1/2019 -> 7/2/2019).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Locate the two author/affiliation placeholders by their current text ---
$authorShape = $null
$affiliationShape = $null

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $t = $sh.TextFrame.TextRange.Text
        if ($t -eq "Priyanka Galla") {
            $authorShape = $sh
        }
        elseif ($t -eq "Northwest Missouri State University") {
            $affiliationShape = $sh
        }
    }
}

# --- "Priyanka Galla" -> "By: Priyanka Galla | Advisor: Dr. Nathan Eloe" (two runs) ---
if ($authorShape -ne $null) {
    $fullText = "By: Priyanka Galla | Advisor: Dr. Nathan Eloe"
    $tr = $authorShape.TextFrame.TextRange
    $tr.Text = $fullText

    $prefix = "By: Priyanka Galla | Advisor: Dr. Nathan "
    $eloeStart = $prefix.Length + 1
    $eloeLen = "Eloe".Length

    # Re-assign the trailing "Eloe" characters to force PowerPoint to split it
    # into its own run (mirrors the authored run-split in the source deck).
    $eloeRange = $tr.Characters($eloeStart, $eloeLen)
    $eloeRange.Text = "Eloe"
}

# --- "Northwest Missouri State University" -> "Email: s534884@nwmissouri.edu" ---
if ($affiliationShape -ne $null) {
    $affiliationShape.TextFrame.TextRange.Text = "Email: s534884@nwmissouri.edu"
}

# --- Re-cache the slide master's date field (7/1/2019 -> 7/2/2019) ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $dsh = $master.Shapes.Item($i)
    if ($dsh.HasTextFrame -and $dsh.TextFrame.TextRange.Text -eq "7/1/2019") {
        $dsh.TextFrame.TextRange.Text = "7/2/2019"
    }
}
